# Auto-generated Excel COM-interop edits matching the target diff
# Updates currentAveragePrice / LevePrice* / LeveProfit* columns (H, I, J, K, L, M, N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets, scraped market data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 575.7
$ws.Range("J17").Value = 588.6896400000001
$ws.Range("L17").Value = 1766.06892
$ws.Range("N17").Value = -2102.06892

$ws.Range("H40").Value = 8173.067
$ws.Range("I40").Value = 3699.2
$ws.Range("K40").Value = 3699.2
$ws.Range("M40").Value = -3524.2

$ws.Range("H113").Value = 2519.4443
$ws.Range("I113").Value = 2475
$ws.Range("J113").Value = 2575
$ws.Range("K113").Value = 2475
$ws.Range("L113").Value = 2575
$ws.Range("M113").Value = 779
$ws.Range("N113").Value = -9083

$ws.Range("H116").Value = 5381.222
$ws.Range("I116").Value = 3845.2727
$ws.Range("K116").Value = 3845.2727
$ws.Range("M116").Value = -403.2727

$ws.Range("H132").Value = 1639.9706
$ws.Range("I132").Value = 1680.0625
$ws.Range("K132").Value = 5040.1875
$ws.Range("M132").Value = -2510.1875

$ws.Range("H133").Value = 77179.914
$ws.Range("J133").Value = 77179.914
$ws.Range("L133").Value = 77179.914
$ws.Range("N133").Value = -87299.914

$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139

$ws.Range("H136").Value = 99995
$ws.Range("J136").Value = 99995
$ws.Range("L136").Value = 99995
$ws.Range("N136").Value = -110195

$ws.Range("H138").Value = 1471.7446
$ws.Range("I138").Value = 884.7586
$ws.Range("J138").Value = 2417.4443
$ws.Range("K138").Value = 2654.2758
$ws.Range("L138").Value = 7252.3329
$ws.Range("M138").Value = 2485.7242
$ws.Range("N138").Value = -17532.3329

$ws.Range("H139").Value = 98406.664
$ws.Range("J139").Value = 98406.664
$ws.Range("L139").Value = 98406.664
$ws.Range("N139").Value = -108686.664

$ws.Range("H140").Value = 80776.664
$ws.Range("J140").Value = 80776.664
$ws.Range("L140").Value = 80776.664
$ws.Range("N140").Value = -91136.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5594.0454
$ws.Range("I32").Value = 2265.0408
$ws.Range("K32").Value = 2265.0408
$ws.Range("M32").Value = -1978.0408

$ws.Range("H74").Value = 2344.1
$ws.Range("I74").Value = 1552.3572
$ws.Range("K74").Value = 1552.3572
$ws.Range("M74").Value = -678.3571999999999

$ws.Range("H77").Value = 2344.1
$ws.Range("I77").Value = 1552.3572
$ws.Range("K77").Value = 7761.786
$ws.Range("M77").Value = -3393.786

$ws.Range("H122").Value = 2555.1428
$ws.Range("I122").Value = 2547.75
$ws.Range("K122").Value = 7643.25
$ws.Range("M122").Value = -5193.25

$ws.Range("H130").Value = 43129.332
$ws.Range("J130").Value = 43129.332
$ws.Range("L130").Value = 43129.332
$ws.Range("N130").Value = -53169.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10488.375
$ws.Range("I86").Value = 7380
$ws.Range("K86").Value = 7380
$ws.Range("M86").Value = -6257

$ws.Range("H89").Value = 10488.375
$ws.Range("I89").Value = 7380
$ws.Range("K89").Value = 36900
$ws.Range("M89").Value = -31284

$ws.Range("H94").Value = 5995
$ws.Range("I94").Value = 5995
$ws.Range("K94").Value = 5995
$ws.Range("M94").Value = -5544

$ws.Range("H99").Value = 37528.355
$ws.Range("I99").Value = 56689.723
$ws.Range("K99").Value = 56689.723
$ws.Range("M99").Value = -55191.723

$ws.Range("H132").Value = 52304.777
$ws.Range("J132").Value = 52304.777
$ws.Range("L132").Value = 52304.777
$ws.Range("N132").Value = -62424.777

$ws.Range("H134").Value = 7333.3125
$ws.Range("I134").Value = 6348.4546
$ws.Range("K134").Value = 19045.3638
$ws.Range("M134").Value = -16510.3638

$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3368.8462
$ws.Range("I31").Value = 2199.2856
$ws.Range("K31").Value = 2199.2856
$ws.Range("M31").Value = -1904.2856

$ws.Range("H34").Value = 3368.8462
$ws.Range("I34").Value = 2199.2856
$ws.Range("K34").Value = 2199.2856
$ws.Range("M34").Value = -1997.2856

$ws.Range("H138").Value = 65242.668
$ws.Range("J138").Value = 65242.668
$ws.Range("L138").Value = 65242.668
$ws.Range("N138").Value = -75522.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 83433
$ws.Range("I26").Value = 152.28572
$ws.Range("J26").Value = 200026
$ws.Range("K26").Value = 456.85716
$ws.Range("L26").Value = 600078
$ws.Range("M26").Value = -168.85716
$ws.Range("N26").Value = -600654

$ws.Range("H63").Value = 2300
$ws.Range("I63").Value = 2300
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 6900
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -6151
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2300
$ws.Range("I66").Value = 2300
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 20700
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16956
$ws.Range("N66").ClearContents()

$ws.Range("H136").Value = 1460.3125
$ws.Range("I136").Value = 1460.3125
$ws.Range("K136").Value = 4380.9375
$ws.Range("M136").Value = 719.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4556496
$ws.Range("J11").Value = 6400000
$ws.Range("L11").Value = 6400000
$ws.Range("N11").Value = -6400278

$ws.Range("H80").Value = 83336216
$ws.Range("J80").Value = 3483.5715
$ws.Range("L80").Value = 3483.5715
$ws.Range("N80").Value = -5479.5715

$ws.Range("H83").Value = 83336216
$ws.Range("J83").Value = 3483.5715
$ws.Range("L83").Value = 17417.8575
$ws.Range("N83").Value = -27401.8575

$ws.Range("H93").Value = 17727
$ws.Range("J93").Value = 17727
$ws.Range("L93").Value = 17727
$ws.Range("N93").Value = -21471

$ws.Range("H109").Value = 28709.154
$ws.Range("J109").Value = 28709.154
$ws.Range("L109").Value = 28709.154
$ws.Range("N109").Value = -30789.154

$ws.Range("H135").Value = 56362.727
$ws.Range("J135").Value = 56362.727
$ws.Range("L135").Value = 56362.727
$ws.Range("N135").Value = -66502.727

$ws.Range("H140").Value = 96372.5
$ws.Range("J140").Value = 97211.42999999999
$ws.Range("L140").Value = 97211.42999999999
$ws.Range("N140").Value = -107571.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3533.8462
$ws.Range("I46").Value = 2593.5
$ws.Range("K46").Value = 2593.5
$ws.Range("M46").Value = -2405.5

$ws.Range("H68").Value = 4542.8335
$ws.Range("I68").Value = 4455.4
$ws.Range("J68").Value = 4980
$ws.Range("K68").Value = 4455.4
$ws.Range("L68").Value = 4980
$ws.Range("M68").Value = -3706.4
$ws.Range("N68").Value = -6478

$ws.Range("H71").Value = 4542.8335
$ws.Range("I71").Value = 4455.4
$ws.Range("J71").Value = 4980
$ws.Range("K71").Value = 22277
$ws.Range("L71").Value = 24900
$ws.Range("M71").Value = -18533
$ws.Range("N71").Value = -32388

$ws.Range("H117").Value = 79225.71000000001
$ws.Range("J117").Value = 84096.664
$ws.Range("L117").Value = 84096.664
$ws.Range("N117").Value = -93274.664

$ws.Range("H132").Value = 3037.889
$ws.Range("I132").Value = 2763.2856
$ws.Range("K132").Value = 8289.856800000001
$ws.Range("M132").Value = -5759.856800000001

$ws.Range("H133").Value = 52288
$ws.Range("J133").Value = 52288
$ws.Range("L133").Value = 52288
$ws.Range("N133").Value = -57348

$ws.Range("H136").Value = 3706
$ws.Range("I136").Value = 6429
$ws.Range("J136").Value = 1890.6666
$ws.Range("K136").Value = 19287
$ws.Range("L136").Value = 19287
$ws.Range("M136").Value = -16737
$ws.Range("N136").Value = -10771.9998
